{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the stack-trace text changes described by the diff:\n//  1) Update 3 lines of frame info (line numbers + GeneratedMethodAccessor index).\n//  2) Replace the Maven Surefire / Tycho / Equinox launcher tail of the\n//     stack trace with the Eclipse JDT JUnit runner tail.\n// The edits live inside a single run's <w:t> in the second body paragraph,\n// so we locate them with Range.search() (unique, exact substrings) and\n// swap them in with Range.insertText(..., Word.InsertLocation.replace).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The stack trace lives in the document's second paragraph (index 1).\nconst stackParagraph = paragraphs.items[1];\n\nconst oldSnippet1 = \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)\\n\\tat sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)\";\nconst newSnippet1 = \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)\\n\\tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)\";\n\nconst oldSnippet2 = \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)\\n\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)\\n\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:161)\\n\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)\\n\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)\\n\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)\\n\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)\\n\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)\\n\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)\\n\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)\\n\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)\\n\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)\\n\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)\\n\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\";\nconst newSnippet2 = \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\";\n\n// --- Edit 1: the three line-number / accessor-index lines ---\nconst found1 = stackParagraph.search(oldSnippet1, { matchCase: true });\nfound1.load(\"items\");\nawait context.sync();\n\nif (found1.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for snippet 1, found \" + found1.items.length);\n}\nfound1.items[0].insertText(newSnippet1, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 2: the long surefire/tycho/equinox tail -> jdt runner tail ---\nconst found2 = stackParagraph.search(oldSnippet2, { matchCase: true });\nfound2.load(\"items\");\nawait context.sync();\n\nif (found2.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for snippet 2, found \" + found2.items.length);\n}\nfound2.items[0].insertText(newSnippet2, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the stack-trace text changes described by the diff:\n#  1) Update 3 lines of frame info (line numbers + GeneratedMethodAccessor index).\n#  2) Replace the Maven Surefire / Tycho / Equinox launcher tail of the\n#     stack trace with the Eclipse JDT JUnit runner tail.\n# Both snippets are unique, literal substrings of the big stack-trace run,\n# so we use Find/Replace (no wildcards) scoped to the whole document content.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: the three line-number / accessor-index lines ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)`n`tat sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)\"\n$find1.Replacement.Text = \"`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)`n`tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)\"\n$find1.Forward = $true\n$find1.Wrap = 1\n$find1.Format = $false\n$find1.MatchCase = $true\n$find1.MatchWildcards = $false\n$replaced1 = $find1.Execute($find1.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\nif (-not $replaced1) {\n    throw \"Snippet 1 not found/replaced\"\n}\n\n# --- Edit 2: the long surefire/tycho/equinox tail -> jdt runner tail ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\"\n$find2.Replacement.Text = \"`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWildcards = $false\n$replaced2 = $find2.Execute($find2.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\nif (-not $replaced2) {\n    throw \"Snippet 2 not found/replaced\"\n}\n\nWrite-Output \"done\"\n"}
